$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new data row (row 4) documenting the "World COVID-19 Events Timeline"
# dataset, mirroring the layout of the existing rows.
# ---------------------------------------------------------------------------

$ws.Range("A4").Value = "World COVID-19 Events Timeline"
$ws.Range("B4").Value = "China Data Lab"
$ws.Range("C4").Value = "https://dataverse.harvard.edu/dataset.xhtml?persistentId=doi:10.7910/DVN/OAM2JK"
$ws.Range("D4").Value = "World"
# E4 (State) intentionally left blank - this dataset has no state scope.
$ws.Range("F4").Value = "day"
$ws.Range("H4").Value = 43955
$ws.Range("I4").Value = "country"
$ws.Range("J4").Value = "Updated to May 3, 2020. Policies and regulations released by the Chinese government, global organizations, western countries, and so on. It is categorized as Chinese News Timeline and Global News Timeline. "
$ws.Range("K4").Value = "Medicine, Health and Life Sciences; Law "

# Turn the Link to Dataset cell into a real hyperlink (adds a row to the
# worksheet's <hyperlinks> collection + an external relationship), then
# restore the cell's look to match the other "Link to Dataset" cells.
$ws.Hyperlinks.Add($ws.Range("C4"), "https://dataverse.harvard.edu/dataset.xhtml?persistentId=doi:10.7910/DVN/OAM2JK")

# Match the formatting of row 3 for the new row (same style per column),
# except the date-start cell (G4), which gets its own plain date format.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("H3:K3").Copy()
$ws.Range("H4:K4").PasteSpecial(-4122)

# G4: Temporal range start - plain short-date format (no special font).
$ws.Range("G4").Value = 43893
$ws.Range("G4").NumberFormat = "mm-dd-yy"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Update the view: scroll right a bit and select L4 (just past the new row).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("L4").Select()
